$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target matrix for B2:K24 (23 rows x 10 columns).
# A new "filter" column of values was inserted as the new column B,
# and the previous B..J columns shifted right into C..K; the former
# column K values are dropped off the right edge (grid stays A1:K24).
$arr = New-Object 'object[,]' 23,10

$arr[0,0] = -0.5351204465965399
$arr[0,1] = -1.66950562919271
$arr[0,2] = -0.2895456687149427
$arr[0,3] = -0.6733414736251095
$arr[0,4] = 0.0115444171491989
$arr[0,5] = -0.6804809672324722
$arr[0,6] = -0.4410326232298434
$arr[0,7] = -0.4043046464706727
$arr[0,8] = 0.4221894279166428
$arr[0,9] = -0.2336087822049224
$arr[1,0] = 0.01136921538350649
$arr[1,1] = -0.3724265895266604
$arr[1,2] = 0.3124593012476481
$arr[1,3] = -0.3795660831340231
$arr[1,4] = -0.1401177391313942
$arr[1,5] = -0.1033897623722235
$arr[1,6] = 0.723104312015092
$arr[1,7] = 0.06730610189352677
$arr[1,8] = -0.5908441378320941
$arr[1,9] = -0.191729189980311
$arr[2,0] = 0.9526635850093691
$arr[2,1] = 0.2606382006276979
$arr[2,2] = 0.5000865446303268
$arr[2,3] = 0.5368145213894975
$arr[2,4] = 1.363308595776813
$arr[2,5] = 0.7075103856552477
$arr[2,6] = 0.0493601459296269
$arr[2,7] = 0.44847509378141
$arr[2,8] = 0.3870385516598326
$arr[2,9] = -0.04712738345727097
$arr[3,0] = 0.4274989419678774
$arr[3,1] = 0.4642269187270481
$arr[3,2] = 1.290720993114364
$arr[3,3] = 0.6349227829927984
$arr[3,4] = -0.02322745673282245
$arr[3,5] = 0.3758874911189606
$arr[3,6] = 0.3144509489973832
$arr[3,7] = -0.1197149861197203
$arr[3,8] = 0.8590864059593566
$arr[3,9] = 0.6319705015114304
$arr[4,0] = 1.235029948750828
$arr[4,1] = 0.5792317386292632
$arr[4,2] = -0.07891850109635756
$arr[4,3] = 0.3201964467554255
$arr[4,4] = 0.2587599046338481
$arr[4,5] = -0.1754060304832554
$arr[4,6] = 0.8033953615958215
$arr[4,7] = 0.5762794571478953
$arr[4,8] = -0.05718027227819322
$arr[4,9] = 0.3702166863774111
$arr[5,0] = -0.4329776612703231
$arr[5,1] = -0.03386271341853997
$arr[5,2] = -0.09529925554011737
$arr[5,3] = -0.529465190657221
$arr[5,4] = 0.449336201421856
$arr[5,5] = 0.2222202969739298
$arr[5,6] = -0.4112394324521587
$arr[5,7] = 0.01615752620344563
$arr[5,8] = -0.1501794284847013
$arr[5,9] = -0.02746841204387546
$arr[6,0] = -0.1856174723396913
$arr[6,1] = -0.6197834074567948
$arr[6,2] = 0.3590179846222821
$arr[6,3] = 0.1319020801743559
$arr[6,4] = -0.5015576492517326
$arr[6,5] = -0.07416069059612829
$arr[6,6] = -0.2404976452842752
$arr[6,7] = -0.1177866288434494
$arr[6,8] = -0.07055289228830908
$arr[6,9] = -0.4671716238107607
$arr[7,0] = 0.4273407187267424
$arr[7,1] = 0.2002248142788162
$arr[7,2] = -0.4332349151472724
$arr[7,3] = -0.005837956491668017
$arr[7,4] = -0.1721749111798149
$arr[7,5] = -0.0494638947389891
$arr[7,6] = -0.002230158183848807
$arr[7,7] = -0.3988488897063004
$arr[7,8] = -0.1447968545825803
$arr[7,9] = 0.07532753529099229
$arr[8,0] = -0.5318964931771777
$arr[8,1] = -0.1044995345215733
$arr[8,2] = -0.2708364892097202
$arr[8,3] = -0.1481254727688944
$arr[8,4] = -0.1008917362137541
$arr[8,5] = -0.4975104677362057
$arr[8,6] = -0.2434584326124856
$arr[8,7] = -0.02333404273891299
$arr[8,8] = -0.2591224913255812
$arr[8,9] = -0.4380952487963659
$arr[9,0] = -0.1938269109680474
$arr[9,1] = -0.07111589452722158
$arr[9,2] = -0.02388215797208129
$arr[9,3] = -0.4205008894945329
$arr[9,4] = -0.1664488543708128
$arr[9,5] = 0.0536755355027598
$arr[9,6] = -0.1821129130839084
$arr[9,7] = -0.3610856705546931
$arr[9,8] = -0.3716462008140141
$arr[9,9] = -0.5793653109721442
$arr[10,0] = 0.1476338940440795
$arr[10,1] = -0.2489848374783721
$arr[10,2] = 0.005067197645347965
$arr[10,3] = 0.2251915875189206
$arr[10,4] = -0.0105968610677476
$arr[10,5] = -0.1895696185385323
$arr[10,6] = -0.2001301487978533
$arr[10,7] = -0.4078492589559834
$arr[10,8] = -0.1906403594810787
$arr[10,9] = 0.1190458097769828
$arr[11,0] = -0.0323979044984018
$arr[11,1] = 0.1877264853751708
$arr[11,2] = -0.04806196321149736
$arr[11,3] = -0.2270347206822821
$arr[11,4] = -0.2375952509416031
$arr[11,5] = -0.4453143610997332
$arr[11,6] = -0.2281054616248284
$arr[11,7] = 0.08158070763323305
$arr[11,8] = -0.07272342619877098
$arr[11,9] = 0.5777029950204122
$arr[12,0] = -0.06996447561954
$arr[12,1] = -0.2489372330903247
$arr[12,2] = -0.2594977633496457
$arr[12,3] = -0.4672168735077758
$arr[12,4] = -0.2500079740328711
$arr[12,5] = 0.05967819522519041
$arr[12,6] = -0.09462593860681362
$arr[12,7] = 0.5558004826123696
$arr[12,8] = 0.3396354339941604
$arr[12,9] = -0.0673936950407959
$arr[13,0] = -0.1713918715036764
$arr[13,1] = -0.3791109816618064
$arr[13,2] = -0.1619020821869017
$arr[13,3] = 0.1477840870711598
$arr[13,4] = -0.006520046760844223
$arr[13,5] = 0.643906374458339
$arr[13,6] = 0.4277413258401298
$arr[13,7] = 0.0207121968051735
$arr[13,8] = 0.5197544139825933
$arr[13,9] = 0.3609055008270807
$arr[14,0] = -0.06818896562035748
$arr[14,1] = 0.241497203637704
$arr[14,2] = 0.08719306980569996
$arr[14,3] = 0.7376194910248832
$arr[14,4] = 0.521454442406674
$arr[14,5] = 0.1144253133717177
$arr[14,6] = 0.6134675305491375
$arr[14,7] = 0.4546186173936249
$arr[14,8] = 0.5107824383638689
$arr[14,9] = 2.760585277975261
$arr[15,0] = 0.2497007499081394
$arr[15,1] = 0.09539661607613537
$arr[15,2] = 0.7458230372953185
$arr[15,3] = 0.5296579886771094
$arr[15,4] = 0.1226288596421531
$arr[15,5] = 0.6216710768195729
$arr[15,6] = 0.4628221636640603
$arr[15,7] = 0.5189859846343043
$arr[15,8] = 2.768788824245696
$arr[15,9] = 10.23793915510299
$arr[16,0] = 0.09280705542466716
$arr[16,1] = 0.7432334766438504
$arr[16,2] = 0.5270684280256412
$arr[16,3] = 0.1200392989906849
$arr[16,4] = 0.6190815161681047
$arr[16,5] = 0.4602326030125921
$arr[16,6] = 0.5163964239828361
$arr[16,7] = 2.766199263594229
$arr[16,8] = 10.23534959445152
$arr[16,9] = -7.935912205685947
$arr[17,0] = 0.7497668092269023
$arr[17,1] = 0.5336017606086931
$arr[17,2] = 0.1265726315737368
$arr[17,3] = 0.6256148487511566
$arr[17,4] = 0.466765935595644
$arr[17,5] = 0.522929756565888
$arr[17,6] = 2.77273259617728
$arr[17,7] = 10.24188292703457
$arr[17,8] = -7.929378873102896
$arr[17,9] = 0.2102926738762539
$arr[18,0] = 0.4223850656296224
$arr[18,1] = 0.01535593659466611
$arr[18,2] = 0.5143981537720859
$arr[18,3] = 0.3555492406165733
$arr[18,4] = 0.4117130615868174
$arr[18,5] = 2.661515901198209
$arr[18,6] = 10.1306662320555
$arr[18,7] = -8.040595568081965
$arr[18,8] = 0.09907597889718328
$arr[18,9] = 2.246646450696576
$arr[19,0] = -0.02639020739223796
$arr[19,1] = 0.4726520097851818
$arr[19,2] = 0.3138030966296693
$arr[19,3] = 0.3699669175999133
$arr[19,4] = 2.619769757211305
$arr[19,5] = 10.0889200880686
$arr[19,6] = -8.082341712068869
$arr[19,7] = 0.05732983491027921
$arr[19,8] = 2.204900306709672
$arr[19,9] = -1.235129679813658
$arr[20,0] = 0.4979670725178967
$arr[20,1] = 0.3391181593623842
$arr[20,2] = 0.3952819803326282
$arr[20,3] = 2.64508481994402
$arr[20,4] = 10.11423515080131
$arr[20,5] = -8.057026649336155
$arr[20,6] = 0.0826448976429941
$arr[20,7] = 2.230215369442386
$arr[20,8] = -1.209814617080943
$arr[20,9] = -1.270988795495144
$arr[21,0] = 0.343156824405298
$arr[21,1] = 0.3993206453755421
$arr[21,2] = 2.649123484986935
$arr[21,3] = 10.11827381584423
$arr[21,4] = -8.052987984293241
$arr[21,5] = 0.08668356268590799
$arr[21,6] = 2.2342540344853
$arr[21,7] = -1.20577595203803
$arr[21,8] = -1.266950130452231
$arr[21,9] = 0.7745058067040239
$arr[22,0] = 0.2804435086845197
$arr[22,1] = 2.530246348295912
$arr[22,2] = 9.999396679153206
$arr[22,3] = -8.171865120984263
$arr[22,4] = -0.03219357400511441
$arr[22,5] = 2.115376897794278
$arr[22,6] = -1.324653088729052
$arr[22,7] = -1.385827267143253
$arr[22,8] = 0.6556286700130015
$arr[22,9] = 0.07303413297936051

$ws.Range("B2:K24").Value = $arr
